$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# New todo item inserted right after the existing block (row 15)
$ws.Range("A15").Value = "Skriv kode referencer"

# Two new todo items appended at the end of the list (rows 24-25)
$ws.Range("A24").Value = "State space"
$ws.Range("A25").Value = "PID"
